$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.606.60'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '2.432.95'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''565.36'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').Value = '''166.71'
$ws.Range('E6').Value = '  +5.79%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').Value = '''0.171'
$ws.Range('E9').Value = '  +8.85%  '
$ws.Range('D10').Value = '2.432.08'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').Value = '''0.335'
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('D13').Value = '''4.70'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('E14').Value = '  +6.24%  '
$ws.Range('D15').Value = '69.200.97'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '2.879.93'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '''24.01'
$ws.Range('E17').Value = '  +5.63%  '
$ws.Range('D18').Value = '2.419.74'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').Value = '''10.84'
$ws.Range('E19').Value = '  +5.81%  '
$ws.Range('D20').Value = '''343.84'
$ws.Range('E20').Value = '  +4.65%  '
$ws.Range('E21').Value = '  +5.73%  '
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('D23').Value = '''1.99'
$ws.Range('E23').Value = '  +7.37%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '''66.13'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  +6.06%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '''8.49'
$ws.Range('E27').Value = '  +6.06%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.558.22'
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''1.01'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').Value = '0.0₃0854'
$ws.Range('E30').Value = '  +8.03%  '
$ws.Range('D31').Value = '''7.38'
$ws.Range('E31').Value = '  +5.28%  '
$ws.Range('D32').Value = '''1.25'
$ws.Range('E32').Value = '  +10.86%  '
$ws.Range('D33').Value = '''455.42'
$ws.Range('E33').Value = '  +9.08%  '
$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  +2.57%  '
$ws.Range('D36').Value = '''159.40'
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '''0.112'
$ws.Range('E37').Value = '  +7.52%  '
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').Value = '''19.11'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  +3.64%  '
$ws.Range('D41').Value = '''0.304'
$ws.Range('E41').Value = '  +4.54%  '
$ws.Range('E42').Value = '  +5.22%  '
$ws.Range('E43').Value = '  +4.86%  '
$ws.Range('D44').Value = '''37.90'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').Value = '''1.09'
$ws.Range('E45').Value = '  +3.24%  '
$ws.Range('D46').Value = '''2.09'
$ws.Range('E46').Value = '  +9.20%  '
$ws.Range('D47').Value = '''136.07'
$ws.Range('E47').Value = '  +5.96%  '
$ws.Range('E48').Value = '  +3.80%  '
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('D50').Value = '''0.491'
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.563'
$ws.Range('E51').Value = '  +1.60%  '
